# Update Name of Algo
# Applies updated KNN-imputed values to column C for specific rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = -13.224
    14 = -12.139
    21 = -12.63
    23 = -12.223
    25 = -12.727
    26 = -13.117
    29 = -12.041
    53 = -11.443
    57 = -13.829
    59 = -13.155
    69 = -10.676
    79 = -12.083
    83 = -13.169
    91 = -10.595
    93 = -11.766
}

foreach ($row in $updates.Keys) {
    $ws.Range("C$row").Value = $updates[$row]
}
